$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.498.00"
$ws.Range("E2").Value = "  +0.90%  "

$ws.Range("D3").Value = "'2.427.47"
$ws.Range("E3").Value = "  +0.62%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'566.09"
$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("D6").Value = "'145.07"
$ws.Range("E6").Value = "  +2.02%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("D11").Value = "'5.30"
$ws.Range("E11").Value = "  +1.15%  "

$ws.Range("D12").Value = "'0.354"
$ws.Range("E12").Value = "  +1.59%  "

$ws.Range("D13").Value = "'26.78"
$ws.Range("E13").Value = "  +4.94%  "

$ws.Range("D14").Value = "'0.0000179"
$ws.Range("E14").Value = "  +4.45%  "

$ws.Range("D15").Value = "'2.793.74"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").Value = "'62.360.71"
$ws.Range("E16").Value = "  +0.85%  "

$ws.Range("D17").Value = "'2.425.71"
$ws.Range("E17").Value = "  +0.69%  "

$ws.Range("D18").Value = "'11.20"
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("D19").Value = "'6.95"
$ws.Range("E19").Value = "  +2.00%  "

$ws.Range("D20").Value = "'323.28"
$ws.Range("E20").Value = "  +0.68%  "

$ws.Range("E21").Value = "  +1.04%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "'66.98"
$ws.Range("E23").Value = "  +1.96%  "

$ws.Range("E24").Value = "  +4.70%  "

$ws.Range("D25").Value = "'595.40"
$ws.Range("E25").Value = "  +5.72%  "

$ws.Range("D26").Value = "'8.55"
$ws.Range("E26").Value = "  -1.24%  "

$ws.Range("D27").Value = "'0.0₃0995"
$ws.Range("E27").Value = "  +7.29%  "

$ws.Range("D28").Value = "'2.543.00"
$ws.Range("E28").Value = "  +0.68%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").Value = "'8.41"
$ws.Range("E30").Value = "  +2.92%  "

$ws.Range("D31").Value = "'1.44"
$ws.Range("E31").Value = "  +3.56%  "

$ws.Range("D32").Value = "'0.144"
$ws.Range("E32").Value = "  -1.98%  "

$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("E34").Value = "  -0.50%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.84"
$ws.Range("E35").Value = "  +1.43%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("D37").Value = "'0.381"
$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("D38").Value = "'18.71"
$ws.Range("E38").Value = "  +1.18%  "

$ws.Range("D39").Value = "'5.35"
$ws.Range("E39").Value = "  -1.61%  "

$ws.Range("D40").Value = "'147.17"
$ws.Range("E40").Value = "  -3.55%  "

$ws.Range("D41").Value = "'1.82"
$ws.Range("E41").Value = "  +1.46%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").Value = "'2.45"
$ws.Range("E43").Value = "  +9.49%  "

$ws.Range("D44").Value = "'147.88"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").Value = "'3.66"
$ws.Range("E45").Value = "  +1.95%  "

$ws.Range("D46").Value = "'0.0534"
$ws.Range("E46").Value = "  +1.04%  "

$ws.Range("D47").Value = "'20.48"
$ws.Range("E47").Value = "  +3.37%  "

$ws.Range("D48").Value = "'0.600"
$ws.Range("E48").Value = "  +1.67%  "

$ws.Range("D49").Value = "'0.0231"
$ws.Range("E49").Value = "  +2.67%  "

$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").Value = "'1.09"
$ws.Range("E51").Value = "  +4.08%  "
